$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) library_content -> library_meta (restructure key/value metadata rows)
# ---------------------------------------------------------------------------
$libMeta = $wb.Worksheets.Item("library_content")

# Clear the old 18-row / 3-column layout entirely before rewriting it.
$libMeta.Range("A1:C18").ClearContents()

$libMeta.Range("A1").Value = "type"
$libMeta.Range("B1").Value = "library"
$libMeta.Range("A2").Value = "urn"
$libMeta.Range("B2").Value = "urn:intuitem:risk:library:tisax-v5.1"
$libMeta.Range("A3").Value = "version"
$libMeta.Range("B3").Value = "'1"
$libMeta.Range("A4").Value = "locale"
$libMeta.Range("B4").Value = "en"
$libMeta.Range("A5").Value = "ref_id"
$libMeta.Range("B5").Value = "TISAX v5.1"
$libMeta.Range("A6").Value = "name"
$libMeta.Range("B6").Value = "Trusted Information Security Assessment Exchange"
$libMeta.Range("A7").Value = "description"
$libMeta.Range("B7").Value = @'
VDA ISA provides the basis for
- a self-assessment to determine the state of information security in an organization (e.g. company)
- audits performed by internal departments (e.g. Internal Audit, Information Security)
- a review in accordance with TISAX (Trusted Information Security Assessment Exchange, http://enx.com/tisax/)
Source: https://portal.enx.com/isa5-en.xlsx
'@
$libMeta.Range("A8").Value = "copyright"
$libMeta.Range("B8").Value = @'
Publisher: VERBAND DER AUTOMOBILINDUSTRIE e. V. (VDA, German Association of the Automotive Industry); Behrenstr. 35; 10117 Berlin; www.vda.de
© 2022 Verband der Automobilindustrie e.V., Berlin
This work has been licensed under Creative Commons Attribution - No Derivative Works 4.0 International Public License. In addition, You are granted the right to distribute derivatives under certain terms.
'@
$libMeta.Range("A9").Value = "provider"
$libMeta.Range("B9").Value = "VDA"
$libMeta.Range("A10").Value = "packager"
$libMeta.Range("B10").Value = "intuitem"

$libMeta.Name = "library_meta"

# ---------------------------------------------------------------------------
# 2) New controls_meta sheet, placed right after library_meta
# ---------------------------------------------------------------------------
$controlsMeta = $wb.Worksheets.Add($null, $libMeta)
$controlsMeta.Name = "controls_meta"
$controlsMeta.Range("A1").Value = "type"
$controlsMeta.Range("B1").Value = "framework"
$controlsMeta.Range("A2").Value = "base_urn"
$controlsMeta.Range("B2").Value = "urn:intuitem:risk:req_node:tisax-v5.1"
$controlsMeta.Range("A3").Value = "urn"
$controlsMeta.Range("B3").Value = "urn:intuitem:risk:framework:tisax-v5.1"
$controlsMeta.Range("A4").Value = "ref_id"
$controlsMeta.Range("B4").Value = "TISAX v5.1"
$controlsMeta.Range("A5").Value = "name"
$controlsMeta.Range("B5").Value = "Trusted Information Security Assessment Exchange"
$controlsMeta.Range("A6").Value = "description"
$controlsMeta.Range("B6").Value = @'
VDA ISA provides the basis for
- a self-assessment to determine the state of information security in an organization (e.g. company)
- audits performed by internal departments (e.g. Internal Audit, Information Security)
- a review in accordance with TISAX (Trusted Information Security Assessment Exchange, http://enx.com/tisax/)
Source: https://portal.enx.com/isa5-en.xlsx
'@
$controlsMeta.Range("A7").Value = "min_score"
$controlsMeta.Range("B7").Value = "'0"
$controlsMeta.Range("A8").Value = "max_score"
$controlsMeta.Range("B8").Value = "'5"
$controlsMeta.Range("A9").Value = "scores_definition"
$controlsMeta.Range("B9").Value = "scores"
$controlsMeta.Range("A10").Value = "implementation_groups_definition"
$controlsMeta.Range("B10").Value = "implementation_groups"

# ---------------------------------------------------------------------------
# 3) controls -> controls_content: drop the empty placeholder cells that
#    littered every row (empty inlineStr cells with no text), then rename.
# ---------------------------------------------------------------------------
$controlsContent = $wb.Worksheets.Item("controls")
$usedRows = $controlsContent.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $controlsContent.Cells.Item($r, $c)
        if ($cell.Value2 -eq "") {
            $cell.ClearContents()
        }
    }
}
$controlsContent.Name = "controls_content"

# ---------------------------------------------------------------------------
# 4) New scores_meta sheet, placed right after controls_content
# ---------------------------------------------------------------------------
$scoresMeta = $wb.Worksheets.Add($null, $controlsContent)
$scoresMeta.Name = "scores_meta"
$scoresMeta.Range("A1").Value = "type"
$scoresMeta.Range("B1").Value = "scores"
$scoresMeta.Range("A2").Value = "name"
$scoresMeta.Range("B2").Value = "scores"

# ---------------------------------------------------------------------------
# 5) scores -> scores_content (content untouched, just renamed)
# ---------------------------------------------------------------------------
$scoresContent = $wb.Worksheets.Item("scores")
$scoresContent.Name = "scores_content"

# ---------------------------------------------------------------------------
# 6) New implementation_groups_meta sheet, placed right after scores_content
# ---------------------------------------------------------------------------
$igMeta = $wb.Worksheets.Add($null, $scoresContent)
$igMeta.Name = "implementation_groups_meta"
$igMeta.Range("A1").Value = "type"
$igMeta.Range("B1").Value = "implementation_groups"
$igMeta.Range("A2").Value = "name"
$igMeta.Range("B2").Value = "implementation_groups"

# ---------------------------------------------------------------------------
# 7) implementation_groups -> implementation_groups_content (renamed only)
# ---------------------------------------------------------------------------
$igContent = $wb.Worksheets.Item("implementation_groups")
$igContent.Name = "implementation_groups_content"

